# "added several new models"
#
# 1) "Medical Devices" sheet: insert a new row for Mimedix (MDXG) right
#    before the existing Hologic (HOLX) row, with a share count of 6.
# 2) "Life Sciences+Diagnostics" sheet: insert a new row for Tempus (TEM)
#    right before the existing Sartorius (SRT GR) row.
# 3) Re-point the GRAL.xlsx hyperlink (it sits on the row that shifted
#    down by one because of the new Tempus row) and refresh view/selection
#    state to match where the edits were made.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Medical Devices: insert Mimedix / MDXG above the Hologic row (row 36)
# ---------------------------------------------------------------------
$wsMed = $wb.Worksheets.Item("Medical Devices")

$wsMed.Rows.Item(36).Insert()
$wsMed.Cells.Item(36, 2).Value = "Mimedix"
$wsMed.Cells.Item(36, 3).Value = "MDXG"
$wsMed.Cells.Item(36, 4).Value = 6
$wsMed.Cells.Item(36, 4).NumberFormat = $wsMed.Cells.Item(4, 4).NumberFormat

# ---------------------------------------------------------------------
# Life Sciences+Diagnostics: insert Tempus / TEM above the Sartorius row
# (row 9), which pushes the Grail hyperlink row from 23 down to 24.
# ---------------------------------------------------------------------
$wsLS = $wb.Worksheets.Item("Life Sciences+Diagnostics")

$wsLS.Rows.Item(9).Insert()
$wsLS.Cells.Item(9, 2).Value = "Tempus"
$wsLS.Cells.Item(9, 3).Value = "TEM"

# Rebuild the two hyperlinks so the Grail one follows its row to B24
# (row-insert shifts the cell content but not the stored hyperlink
# range, and the GRAL.xlsx link used to live on B23).
$wsLS.Range("A1").Hyperlinks.Delete()
$wsLS.Hyperlinks.Add($wsLS.Range("B3"), "TMO.xlsx")
$wsLS.Hyperlinks.Add($wsLS.Range("B24"), "GRAL.xlsx")

# ---------------------------------------------------------------------
# View/selection state: Medical Devices becomes the active tab, with the
# new row selected; Life Sciences+Diagnostics keeps its own selection;
# Animal's last selection also moved.
# ---------------------------------------------------------------------
$wsLS.Range("A10").Select()

$wsAnimal = $wb.Worksheets.Item("Animal")
$wsAnimal.Range("C3").Select()

$wsMed.Activate()
$wsMed.Range("D36").Select()

Write-Output "done"
